$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "time_taken" header in F1, matching the style of the other headers (E1).
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

# Populate F2:F89 with the recorded time_taken values.
$timestamps = @(
    "2021-10-05 13:42:35.957392",
    "2021-10-05 13:42:35.957406",
    "2021-10-05 13:42:35.957410",
    "2021-10-05 13:42:35.957413",
    "2021-10-05 13:42:35.957416",
    "2021-10-05 13:42:35.957419",
    "2021-10-05 13:42:35.957422",
    "2021-10-05 13:42:35.957425",
    "2021-10-05 13:42:35.957428",
    "2021-10-05 13:42:35.957431",
    "2021-10-05 13:42:35.957434",
    "2021-10-05 13:42:35.957437",
    "2021-10-05 13:42:35.957440",
    "2021-10-05 13:42:35.957443",
    "2021-10-05 13:42:35.957446",
    "2021-10-05 13:42:35.957449",
    "2021-10-05 13:42:35.957452",
    "2021-10-05 13:42:35.957455",
    "2021-10-05 13:42:35.957459",
    "2021-10-05 13:42:35.957461",
    "2021-10-05 13:42:35.957464",
    "2021-10-05 13:42:35.957467",
    "2021-10-05 13:42:35.957470",
    "2021-10-05 13:42:35.957473",
    "2021-10-05 13:42:35.957476",
    "2021-10-05 13:42:35.957479",
    "2021-10-05 13:42:35.957482",
    "2021-10-05 13:42:35.957485",
    "2021-10-05 13:42:35.957488",
    "2021-10-05 13:42:35.957491",
    "2021-10-05 13:42:35.957494",
    "2021-10-05 13:42:35.957497",
    "2021-10-05 13:42:35.957501",
    "2021-10-05 13:42:35.957504",
    "2021-10-05 13:42:35.957507",
    "2021-10-05 13:42:35.957510",
    "2021-10-05 13:42:35.957513",
    "2021-10-05 13:42:35.957516",
    "2021-10-05 13:42:35.957519",
    "2021-10-05 13:42:35.957522",
    "2021-10-05 13:42:35.957525",
    "2021-10-05 13:42:35.957529",
    "2021-10-05 13:42:35.957532",
    "2021-10-05 13:42:35.957535",
    "2021-10-05 13:42:35.957538",
    "2021-10-05 13:42:35.957541",
    "2021-10-05 13:42:35.957544",
    "2021-10-05 13:42:35.957546",
    "2021-10-05 13:42:35.957549",
    "2021-10-05 13:42:35.957552",
    "2021-10-05 13:42:35.957555",
    "2021-10-05 13:42:35.957558",
    "2021-10-05 13:42:35.957562",
    "2021-10-05 13:42:35.957565",
    "2021-10-05 13:42:35.957568",
    "2021-10-05 13:42:35.957571",
    "2021-10-05 13:42:35.957574",
    "2021-10-05 13:42:35.957577",
    "2021-10-05 13:42:35.957580",
    "2021-10-05 13:42:35.957583",
    "2021-10-05 13:42:35.957586",
    "2021-10-05 13:42:35.957589",
    "2021-10-05 13:42:35.957592",
    "2021-10-05 13:42:35.957595",
    "2021-10-05 13:42:35.957600",
    "2021-10-05 13:42:35.957603",
    "2021-10-05 13:42:35.957606",
    "2021-10-05 13:42:35.957609",
    "2021-10-05 13:42:35.957612",
    "2021-10-05 13:42:35.957615",
    "2021-10-05 13:42:35.957618",
    "2021-10-05 13:42:35.957621",
    "2021-10-05 13:42:35.957625",
    "2021-10-05 13:42:35.957628",
    "2021-10-05 13:42:35.957631",
    "2021-10-05 13:42:35.957634",
    "2021-10-05 13:42:35.957639",
    "2021-10-05 13:42:35.957643",
    "2021-10-05 13:42:35.957646",
    "2021-10-05 13:42:35.957649",
    "2021-10-05 13:42:35.957652",
    "2021-10-05 13:42:35.957655",
    "2021-10-05 13:42:35.957658",
    "2021-10-05 13:42:35.957661",
    "2021-10-05 13:42:35.957664",
    "2021-10-05 13:42:35.957667",
    "2021-10-05 13:42:35.957670",
    "2021-10-05 13:42:35.957673"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
